# Update cryptos list prices and 1h volume percentages (Price column D,
# Volume(1h) column E) with freshly scraped values. D-column cells store
# the price as TEXT (matching the sheet's original t="inlineStr" cells),
# so numeric-looking prices are written with a leading quote to keep
# Excel from auto-converting them to numbers, then restored to the
# workbook's default "Normal" style so no visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.211.14"
$ws.Range("E2").Value = "  +5.44%  "
$ws.Range("D3").Value = "2.743.32"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'578.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'154.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.48%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").Value = "2.757.73"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("D12").Value = "'0.388"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("E13").Value = "  +3.70%  "
$ws.Range("D14").Value = "3.237.11"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").Value = "'26.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "64.099.00"
$ws.Range("E16").Value = "  +5.30%  "
$ws.Range("E17").Value = "  +6.49%  "
$ws.Range("D18").Value = "2.751.54"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").Value = "'11.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("D21").Value = "'360.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").Value = "'6.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'0.531"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "'66.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("D26").Value = "'0.171"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.62%  "
$ws.Range("D27").Value = "'8.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.93%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "0.0₃0916"
$ws.Range("E29").Value = "  +12.56%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "'7.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.65%  "
$ws.Range("D32").Value = "'1.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +19.44%  "
$ws.Range("D33").Value = "'172.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'20.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").Value = "'4.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.75%  "
$ws.Range("E37").Value = "  +8.36%  "
$ws.Range("D38").Value = "'1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.68%  "
$ws.Range("D39").Value = "'1.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +16.12%  "
$ws.Range("D40").Value = "'344.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("D42").Value = "'39.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").Value = "'5.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.81%  "
$ws.Range("D44").Value = "'21.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.35%  "
$ws.Range("D45").Value = "'21.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.00%  "
$ws.Range("D46").Value = "'0.0589"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.98%  "
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("D48").Value = "'138.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
